$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation for "Macroferia Regional de Talca - Berenjena"
# is inserted at row 45, pushing the existing rows 45-65 down to rows 46-66.
$ws.Rows.Item(45).Insert()

# Populate the newly inserted row 45 with the new observation's data.
$ws.Cells.Item(45, 1).Value  = 5
$ws.Cells.Item(45, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(45, 3).Value  = "Maule"
$ws.Cells.Item(45, 4).Value  = 44466
$ws.Cells.Item(45, 5).Value  = 7
$ws.Cells.Item(45, 6).Value  = 100112001
$ws.Cells.Item(45, 7).Value  = "Berenjena"
$ws.Cells.Item(45, 8).Value  = "Sin especificar"
$ws.Cells.Item(45, 9).Value  = "Primera"
$ws.Cells.Item(45, 10).Value = 300
$ws.Cells.Item(45, 11).Value = 7000
$ws.Cells.Item(45, 12).Value = 7000
$ws.Cells.Item(45, 13).Value = 7000
$ws.Cells.Item(45, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(45, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(45, 16).Value = 140
$ws.Cells.Item(45, 17).Value = 50
$ws.Cells.Item(45, 18).Value = "Hortaliza"

# Make sure the date column keeps its date/time number format.
$ws.Cells.Item(45, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
